$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 26156
$ws.Range("E2").Value = 79
$ws.Range("F2").Value = 79
$ws.Range("G2").Value = -390
$ws.Range("H2").Value = -1182
$ws.Range("I2").Value = -1188
$ws.Range("J2").Value = 6
$ws.Range("K2").Value = 36017
$ws.Range("L2").Value = 19587
$ws.Range("M2").Value = 16431
$ws.Range("N2").Value = 16412
$ws.Range("O2").Value = 18
$ws.Range("P2").Value = 2656
$ws.Range("Q2").Value = 114
$ws.Range("R2").Value = -1029
$ws.Range("S2").Value = -897
$ws.Range("T2").Value = 1137
$ws.Range("U2").Value = -1023
$ws.Range("V2").Value = 6102
$ws.Range("W2").Value = 0.3
$ws.Range("X2").Value = -4.52
$ws.Range("Y2").Value = -7.11
$ws.Range("Z2").Value = -3.36
$ws.Range("AA2").Value = 119.21
$ws.Range("AB2").Value = 450.15
$ws.Range("AC2").Value = -2236
$ws.Range("AD2").Value = -10.66
$ws.Range("AE2").Value = 30891
$ws.Range("AF2").Value = 0.77
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 53130000
$ws.Range("D3").Value = 26134
$ws.Range("E3").Value = -596
$ws.Range("F3").Value = -596
$ws.Range("G3").Value = -191
$ws.Range("H3").Value = 62
$ws.Range("I3").Value = 57
$ws.Range("J3").Value = 5
$ws.Range("K3").Value = 41079
$ws.Range("L3").Value = 18628
$ws.Range("M3").Value = 22451
$ws.Range("N3").Value = 22434
$ws.Range("O3").Value = 18
$ws.Range("P3").Value = 2656
$ws.Range("Q3").Value = -375
$ws.Range("R3").Value = 1776
$ws.Range("S3").Value = -918
$ws.Range("T3").Value = 419
$ws.Range("U3").Value = -794
$ws.Range("V3").Value = 5216
$ws.Range("W3").Value = -2.28
$ws.Range("X3").Value = 0.24
$ws.Range("Y3").Value = 0.3
$ws.Range("Z3").Value = 0.16
$ws.Range("AA3").Value = 82.97
$ws.Range("AB3").Value = 453.3
$ws.Range("AC3").Value = 108
$ws.Range("AD3").Value = 329.21
$ws.Range("AE3").Value = 42224
$ws.Range("AF3").Value = 0.84
$ws.Range("AG3").Value = 300
$ws.Range("AH3").Value = 0.84
$ws.Range("AI3").Value = 277.81
$ws.Range("AJ3").Value = 53130000
$ws.Range("D4").Value = 35189
$ws.Range("E4").Value = 1507
$ws.Range("F4").Value = 1507
$ws.Range("G4").Value = 4185
$ws.Range("H4").Value = 3459
$ws.Range("I4").Value = 3452
$ws.Range("J4").Value = 6
$ws.Range("K4").Value = 56524
$ws.Range("L4").Value = 33117
$ws.Range("M4").Value = 23408
$ws.Range("N4").Value = 23390
$ws.Range("O4").Value = 17
$ws.Range("P4").Value = 2656
$ws.Range("Q4").Value = 1961
$ws.Range("R4").Value = -7183
$ws.Range("S4").Value = 6379
$ws.Range("T4").Value = 891
$ws.Range("U4").Value = 1069
$ws.Range("V4").Value = 11858
$ws.Range("W4").Value = 4.28
$ws.Range("X4").Value = 9.83
$ws.Range("Y4").Value = 15.07
$ws.Range("Z4").Value = 7.09
$ws.Range("AA4").Value = 141.48
$ws.Range("AB4").Value = 579.14
$ws.Range("AC4").Value = 6498
$ws.Range("AD4").Value = 6.69
$ws.Range("AE4").Value = 44025
$ws.Range("AF4").Value = 0.99
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 53130000
$ws.Range("D5").Value = 42155
$ws.Range("E5").Value = 829
$ws.Range("F5").Value = 829
$ws.Range("G5").Value = 178
$ws.Range("H5").Value = -477
$ws.Range("I5").Value = -482
$ws.Range("J5").Value = 5
$ws.Range("K5").Value = 57476
$ws.Range("L5").Value = 35666
$ws.Range("M5").Value = 21810
$ws.Range("N5").Value = 21794
$ws.Range("O5").Value = 16
$ws.Range("P5").Value = 2656
$ws.Range("Q5").Value = -97
$ws.Range("R5").Value = 726
$ws.Range("S5").Value = 1404
$ws.Range("T5").Value = 1416
$ws.Range("U5").Value = -1514
$ws.Range("V5").Value = 13447
$ws.Range("W5").Value = 1.97
$ws.Range("X5").Value = -1.13
$ws.Range("Y5").Value = -2.14
$ws.Range("Z5").Value = -0.84
$ws.Range("AA5").Value = 163.53
$ws.Range("AB5").Value = 556.35
$ws.Range("AC5").Value = -913
$ws.Range("AD5").Value = -38.92
$ws.Range("AE5").Value = 41433
$ws.Range("AF5").Value = 0.86
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 52600000
$ws.Range("D6").Value = 44532
$ws.Range("E6").Value = 532
$ws.Range("F6").Value = 532
$ws.Range("G6").Value = 535
$ws.Range("H6").Value = 493
$ws.Range("I6").Value = 323
$ws.Range("K6").Value = 73967
$ws.Range("L6").Value = 47609
$ws.Range("M6").Value = 26358
$ws.Range("N6").Value = 23696
$ws.Range("P6").Value = 2656
$ws.Range("Q6").Value = 563
$ws.Range("R6").Value = -2131
$ws.Range("S6").Value = 3327
$ws.Range("T6").Value = 1725
$ws.Range("U6").Value = -1161
$ws.Range("V6").Value = 17341
$ws.Range("W6").Value = 1.19
$ws.Range("X6").Value = 1.11
$ws.Range("Y6").Value = 1.42
$ws.Range("Z6").Value = 0.75
$ws.Range("AA6").Value = 180.63
$ws.Range("AB6").Value = 667.34
$ws.Range("AC6").Value = 619
$ws.Range("AD6").Value = 52.55
$ws.Range("AE6").Value = 45535
$ws.Range("AF6").Value = 0.71
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 52040000
$ws.Range("D7").Value = 53741
$ws.Range("E7").Value = 1749
$ws.Range("G7").Value = 1753
$ws.Range("H7").Value = 1674
$ws.Range("I7").Value = 1375
$ws.Range("K7").Value = 81552
$ws.Range("L7").Value = 53667
$ws.Range("M7").Value = 27885
$ws.Range("N7").Value = 24998
$ws.Range("P7").Value = 2651
$ws.Range("Q7").Value = 3603
$ws.Range("R7").Value = -5297
$ws.Range("S7").Value = 1293
$ws.Range("T7").Value = 1783
$ws.Range("U7").Value = 1260
$ws.Range("W7").Value = 3.25
$ws.Range("X7").Value = 3.11
$ws.Range("Y7").Value = 5.65
$ws.Range("Z7").Value = 2.15
$ws.Range("AA7").Value = 192.46
$ws.Range("AC7").Value = 2658
$ws.Range("AD7").Value = 12.23
$ws.Range("AE7").Value = 48483
$ws.Range("AF7").Value = 0.67
$ws.Range("AG7").Value = 20
$ws.Range("AH7").Value = 0.06
$ws.Range("AI7").Value = 0.75
$ws.Range("D8").Value = 59179
$ws.Range("E8").Value = 2337
$ws.Range("G8").Value = 2078
$ws.Range("H8").Value = 1642
$ws.Range("I8").Value = 1314
$ws.Range("K8").Value = 84971
$ws.Range("L8").Value = 55522
$ws.Range("M8").Value = 29449
$ws.Range("N8").Value = 26323
$ws.Range("P8").Value = 2651
$ws.Range("Q8").Value = 2772
$ws.Range("R8").Value = -2774
$ws.Range("S8").Value = -57
$ws.Range("T8").Value = 1762
$ws.Range("U8").Value = 927
$ws.Range("W8").Value = 3.95
$ws.Range("X8").Value = 2.77
$ws.Range("Y8").Value = 5.12
$ws.Range("Z8").Value = 1.97
$ws.Range("AA8").Value = 188.54
$ws.Range("AC8").Value = 2548
$ws.Range("AD8").Value = 12.75
$ws.Range("AE8").Value = 51054
$ws.Range("AF8").Value = 0.64
$ws.Range("AG8").Value = 55
$ws.Range("AH8").Value = 0.17
$ws.Range("AI8").Value = 2.16
$ws.Range("D9").Value = 62361
$ws.Range("E9").Value = 2739
$ws.Range("G9").Value = 2530
$ws.Range("H9").Value = 1997
$ws.Range("I9").Value = 1590
$ws.Range("K9").Value = 87748
$ws.Range("L9").Value = 56358
$ws.Range("M9").Value = 31390
$ws.Range("N9").Value = 27845
$ws.Range("P9").Value = 2650
$ws.Range("Q9").Value = 3638
$ws.Range("R9").Value = -2746
$ws.Range("S9").Value = -289
$ws.Range("T9").Value = 1755
$ws.Range("U9").Value = 1742
$ws.Range("W9").Value = 4.39
$ws.Range("X9").Value = 3.2
$ws.Range("Y9").Value = 5.87
$ws.Range("Z9").Value = 2.31
$ws.Range("AA9").Value = 179.54
$ws.Range("AC9").Value = 3084
$ws.Range("AD9").Value = 10.54
$ws.Range("AE9").Value = 54005
$ws.Range("AF9").Value = 0.6
$ws.Range("AG9").Value = 67
$ws.Range("AH9").Value = 0.21
$ws.Range("AI9").Value = 2.16
